# Add a new worksheet "header_error" after the last sheet ("lookup"),
# populate it with a small table that has a duplicated header label
# (COLUMN_A used twice), then make "general" the active/selected sheet
# again (it was "lookup" before).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "header_error"

# Header row (row 2) - note COLUMN_A repeats in B2 and D2.
$newSheet.Range("B2").Value = "COLUMN_A"
$newSheet.Range("C2").Value = "COLUMN_B"
$newSheet.Range("D2").Value = "COLUMN_A"

# Column B data first ...
$newSheet.Range("B3").Value = 10
$newSheet.Range("B4").Value = 20

# ... then column C ...
$newSheet.Range("C3").Value = "hey"
$newSheet.Range("C4").Value = "you"

# ... then column D, to match the original shared-string insertion order.
$newSheet.Range("D3").Value = "out"
$newSheet.Range("D4").Value = "there"

[void]$newSheet.Range("F6").Select()

# Restore "general" as the active sheet (it was "lookup" previously).
$general = $wb.Worksheets.Item("general")
$general.Activate()
